$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# "Advance TestNG" row: Status column (E) for the vtiger row moves from the
# stray numeric placeholder "123" to the real expected result "pass".
$ws.Range("E2").Value = "pass"

# Author left the selection sitting on the edited cell before saving.
$ws.Range("E2").Select()

# Author also zoomed the sheet view in to 160% before saving.
$excel.ActiveWindow.Zoom = 160
